$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing data rows (rows 4 and 5) - table shrinks to 2 data rows
$ws.Rows("4:5").Delete()

# Remove the "Nodes" column (old column D) - "Hidden Layers" (C) and "Nodes" (D)
# are merged into a single "HL Nodes" column; deleting D shifts E->D, F->E
$ws.Columns("D:D").Delete()

# Set the new "(5,3)" node-config values before the "HL Nodes" header text so
# shared-string interning order matches the author's edit order
$ws.Range("C2").Value = "(5,3)"
$ws.Range("C3").Value = "(5,3)"

# Update header row
$ws.Range("C1").Value = "HL Nodes"
$ws.Range("D1").Value = "Error Train"
$ws.Range("E1").Value = "Error Test"

# Update row 2 (first network configuration)
$ws.Range("A2").Value = "stochastic gradient Descent"
$ws.Range("B2").Value = "ReLU"
$ws.Range("D2").Value = 0.326
$ws.Range("E2").Value = 0.341

# Update row 3 (second network configuration)
$ws.Range("A3").Value = "stochastic gradient Descent"
$ws.Range("B3").Value = "Logistic"
$ws.Range("D3").Value = 0.431
$ws.Range("E3").Value = 0.447

# Update active selection to match the author's cursor position
$ws.Range("C8").Select()
